$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format ("@") on all touched Price/Volume cells so that
# numeric-looking values (e.g. "310.24", "-2.76%") are stored as literal
# text, matching the workbook convention used throughout this sheet
# (every data cell is an inline/shared string, never a numeric cell).
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "E14", "D15", "E15", "E16", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "D26", "E26", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume 1h (E) figures.
$ws.Range("D2").Value = "310.24"
$ws.Range("E2").Value = "-2.76%"
$ws.Range("D3").Value = "52.03"
$ws.Range("E3").Value = "6.95%"
$ws.Range("D4").Value = "5.117"
$ws.Range("E4").Value = "-2.74%"
$ws.Range("D5").Value = "0.07793"
$ws.Range("E5").Value = "-2.82%"
$ws.Range("D6").Value = "4.500"
$ws.Range("E6").Value = "-1.85%"
$ws.Range("D7").Value = "1.359"
$ws.Range("E7").Value = "-4.90%"
$ws.Range("D8").Value = "1.584"
$ws.Range("E8").Value = "-3.83%"
$ws.Range("D9").Value = "0.1217"
$ws.Range("E9").Value = "-4.80%"
$ws.Range("D10").Value = "0.2011"
$ws.Range("E10").Value = "3.47%"
$ws.Range("D11").Value = "0.04731"
$ws.Range("E11").Value = "2.88%"
$ws.Range("D12").Value = "0.09480"
$ws.Range("E12").Value = "1.70%"
$ws.Range("D13").Value = "0.1046"
$ws.Range("E13").Value = "0.20%"
$ws.Range("E14").Value = "-4.34%"
$ws.Range("D15").Value = "0.005795"
$ws.Range("E15").Value = "-0.67%"
$ws.Range("E16").Value = "2,017.75%"
$ws.Range("E17").Value = "0.10%"
$ws.Range("D18").Value = "2.439"
$ws.Range("E18").Value = "-0.12%"
$ws.Range("D19").Value = "0.3472"
$ws.Range("E19").Value = "1.51%"
$ws.Range("D20").Value = "8.002"
$ws.Range("E20").Value = "-1.91%"
$ws.Range("D21").Value = "0.1369"
$ws.Range("E21").Value = "-2.12%"
$ws.Range("D22").Value = "0.3094"
$ws.Range("E22").Value = "-0.02%"
$ws.Range("D23").Value = "0.04170"
$ws.Range("E23").Value = "0.09%"
$ws.Range("D24").Value = "0.001263"
$ws.Range("E24").Value = "-3.31%"
$ws.Range("E25").Value = "-7.49%"
$ws.Range("D26").Value = "0.0001351"
$ws.Range("E26").Value = "0.04%"
$ws.Range("E38").Value = "-3.92%"
$ws.Range("D39").Value = "0.05896"
$ws.Range("E39").Value = "3.42%"
$ws.Range("D40").Value = "0.01073"
$ws.Range("E40").Value = "69.90%"
$ws.Range("D41").Value = "0.008078"
$ws.Range("E41").Value = "0.92%"
$ws.Range("D42").Value = "0.1421"
$ws.Range("E42").Value = "-1.32%"
$ws.Range("D43").Value = "0.008247"
$ws.Range("E43").Value = "7.41%"
$ws.Range("D44").Value = "0.008438"
$ws.Range("E44").Value = "6.91%"
$ws.Range("D45").Value = "0.3121"
$ws.Range("E45").Value = "-10.63%"
$ws.Range("D46").Value = "0.00007363"
$ws.Range("E46").Value = "6.72%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.12%"
$ws.Range("D48").Value = "0.05730"
$ws.Range("E48").Value = "4.45%"
$ws.Range("D49").Value = "0.002623"
$ws.Range("E49").Value = "-34.34%"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "0.12%"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "0.12%"
